$d = $word.ActiveDocument

# --- Body text: gendoc template output path ---
# "...\Gendoc\pureEthernetStructure_2.0.0-tsi.d.t+gendoc.${date}.${time}docx'"
# becomes
# "...\Gendoc\PureEthernetStructure_2.0.0-tsi.d.t+gendoc.1.docx'"
$d.Content.Find.Execute("Gendoc\pureEthernetStructure", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Gendoc\PureEthernetStructure", 2)
$d.Content.Find.Execute("gendoc.`${date}.`${time}docx", $true, $false, $false, $false, $false, `
    $true, 1, $false, "gendoc.1.docx", 2)

# --- Header text: version / build tag ---
# "PureEthernetStructure<TAB>2.0.0-tsi.d.t+gendoc.n" becomes
# "PureEthernetStructure<TAB>2.0.0-tsi.d.t+gendoc.1"
$hdr = $d.Sections(1).Headers(1)
$hdr.Range.Find.Execute("gendoc.n", $true, $false, $false, $false, $false, `
    $true, 1, $false, "gendoc.1", 2)
